# Insert a new weekly price-report row for "Terminal La Palmera de La Serena - Apio"
# (Hortaliza), matching the new source date 45147 (2023-08-09), shifting all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 658; everything from the old row 658 downward
# shifts to row+1 automatically.
$ws.Rows.Item(658).Insert()

# Seed the new row 658 with the same record that is now sitting in row 659
# (a duplicate of the data that used to be at row 658), then overwrite the
# two cells that actually differ for the new week: Fecha (D) and Volumen (J).
$ws.Range("A659:R659").Copy()
$ws.Range("A658").PasteSpecial()

$ws.Cells.Item(658, 4).Value = 45147
$ws.Cells.Item(658, 10).Value = 1600

Write-Host "Inserted row 658 -> Fecha:" $ws.Cells.Item(658, 4).Value2 "Volumen:" $ws.Cells.Item(658, 10).Value2
Write-Host "New dimension used rows:" $ws.UsedRange.Rows.Count
